# Lead_ExistingOpportunity.xlsx - add 5 new "Filter" test-case sheets
# (TC11-TC15) cloned from the existing ExistingOpportunityFilter_TC10
# sheet, fix up a couple of stray spaces in existing shared strings, and
# leave the new TC15 sheet active/selected (matches the author's manual
# edit captured in the commit "New added script and Updated Commit =
# 29/09/2020").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Trim stray whitespace on two pre-existing labels.
# ---------------------------------------------------------------------
$wsTC05 = $wb.Worksheets.Item("ExistingOpportunity_TC05")
$wsTC05.Range("D2").Value = "Test User"

$wsTC10 = $wb.Worksheets.Item("ExistingOpportunityFilter_TC10")
$wsTC10.Range("D2").Value = "Test MIRketa1"

# ---------------------------------------------------------------------
# 2. Build the five new sheets by cloning the TC10 template (same
#    columns/styles/header row) and appending them after the last sheet.
# ---------------------------------------------------------------------
# NOTE: this PowerShell engine does not reliably bind named (-Param)
# arguments to custom functions, so every call below uses positional
# arguments instead: SourceName, NewName, D2, E2, F2, F2Styled, H2.
function Add-FilterSheet {
    param(
        [string]$SourceName,
        [string]$NewName,
        [string]$D2,
        [string]$E2,
        [string]$F2,
        [bool]$F2Styled,
        [string]$H2
    )

    $source = $wb.Worksheets.Item($SourceName)
    $source.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
    $ns = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ns.Name = $NewName

    # All five new sheets use A2 = 3000 (TC10's template carries 2000).
    $ns.Range("A2").Value = 3000

    $ns.Range("D2").Value = $D2
    $ns.Range("E2").Value = $E2
    $ns.Range("H2").Value = $H2

    if ($F2Styled) {
        $ns.Range("F2").Value = $F2
    } else {
        # TC13/TC14/TC15 lose the "Cold" number-format styling that the
        # TC10 template carries on F2 - clear it before writing "Hot".
        $ns.Range("F2").ClearFormats()
        $ns.Range("F2").Value = $F2
    }

    return $ns
}

$tc11 = Add-FilterSheet "ExistingOpportunityFilter_TC10" "ExistingOpportunityFilter_TC11" "Test Admin" "Lead_EO_TC11" "Cold" $true "Web"

$tc12 = Add-FilterSheet "ExistingOpportunityFilter_TC11" "ExistingOpportunityFilter_TC12" "Test Host" "Lead_EO_TC12" "Cold" $true "Web"

$tc13 = Add-FilterSheet "ExistingOpportunityFilter_TC12" "ExistingOpportunityFilter_TC13" "Test Provider" "Lead_EO_TC13" "Hot" $false "Web"

$tc14 = Add-FilterSheet "ExistingOpportunityFilter_TC13" "ExistingOpportunityFilter_TC14" "Test Advisor" "Lead_EO_TC14" "Hot" $false "Phone Inquiry"

$tc15 = Add-FilterSheet "ExistingOpportunityFilter_TC14" "ExistingOpportunityFilter_TC15" "Test Banking" "Lead_EO_TC15" "Hot" $false "Web"

# ---------------------------------------------------------------------
# 3. Restore the selections the author left on screen: TC10 parked on
#    D2, TC11 parked on E2, and TC15 (the newest/last sheet) active with
#    A2 selected.
# ---------------------------------------------------------------------
$wsTC10.Select()
$wsTC10.Range("D2").Select()

$tc11.Select()
$tc11.Range("E2").Select()

$tc15.Select()
$tc15.Range("A2").Select()

Write-Host "Added TC11-TC15 and refreshed TC05/TC10 labels"
